$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.149.62"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "3.096.08"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").Value = "3.092.47"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.476"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "3.616.23"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "67.081.29"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.091.84"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "489.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.78%  "
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").Value = "0.0₃0945"
$ws.Range("E34").Value = "  -5.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.966"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.35%  "
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.306"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("D43").Value = "2.790.17"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "381.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.95%  "
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "134.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("E51").Value = "  -2.28%  "
